$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (precision)
$ws.Range("C2").Value = 0.25
$ws.Range("F2").Value = 0.1666666666666667
$ws.Range("U2").Value = 1
$ws.Range("W2").Value = 0.5

# Row 3 (recall)
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("W3").Value = 1

# Row 4 (f1-score) - previously empty inline strings, now numeric
$ws.Range("C4").Value = 0.4
$ws.Range("F4").Value = 0.2857142857142857
$ws.Range("U4").Value = 1
$ws.Range("W4").Value = 0.6666666666666666

# Row 5 (f2-score) - previously empty inline strings, now numeric
$ws.Range("C5").Value = 0.625
$ws.Range("F5").Value = 0.5
$ws.Range("U5").Value = 1
$ws.Range("W5").Value = 0.8333333333333334

# Row 6 (NDCG)
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("W6").Value = 1

# Row 7 (M1) - boolean
$ws.Range("C7").Value = $true
$ws.Range("F7").Value = $true
$ws.Range("U7").Value = $true
$ws.Range("W7").Value = $true

# Row 8 (M3) - boolean
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("U8").Value = $true
$ws.Range("W8").Value = $true

# Row 9 (M5) - boolean
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("U9").Value = $true
$ws.Range("W9").Value = $true

# Row 10 (position) - previously empty inline strings, now numeric
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("W10").Value = 1
